$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A from 19 to 20 (stored OOXML width).
# Excel's ColumnWidth property (character units) is offset from the raw
# stored "width" attribute by the sheet's pixel-rounding; 19.16 round-trips
# to a stored width of exactly 20 for this workbook's default font.
$ws.Columns.Item(1).ColumnWidth = 19.16

# Header change
$ws.Range("A1").Value = "Suchergebnis"

# Find the last used row in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    # NOTE: -eq in this engine is case-insensitive, so we use the
    # case-sensitive .Equals() method to avoid mismatching cells that
    # only differ by letter case (e.g. "m200-222448" vs "M043-17470").
    if ($val.Equals("1_Suchfeld ist leer")) {
        $cell.Value = "1_Das zu durchsuchende Feld ist leer"
    }
    elseif ($val.Equals("1_Nichts gefunden")) {
        $cell.Value = "1_Keine Reparaturnummer gefunden"
    }
    elseif ($val.Equals(" 107-13068")) {
        $cell.Value = "1_Keine Reparaturnummer gefunden"
    }
    elseif ($val.Equals("M043-17470")) {
        $cell.Value = "1_Keine Reparaturnummer gefunden"
    }
    elseif ($val.Equals("m200-222448")) {
        $cell.Value = "M200-222448"
    }
    elseif ($val.Equals("s021-12695")) {
        $cell.Value = "S021-12695"
    }
}
